$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each target cell currently holds a text value (e.g. "38.734.53" or "  -4.77%  ").
# Assigning a plain numeric-looking string via .Value would let Excel auto-convert
# it into a real number, changing the cell type. Prefixing the text with a leading
# apostrophe forces it to be stored as text (same as a user typing it in Excel).
# Resetting .Style to "Normal" afterwards removes the quote-prefix "Text" number
# format that the apostrophe entry applies, restoring the cell to its original,
# unstyled state while keeping the exact text value.

$c = $ws.Range("D2")
$c.Value = "'" + '38.725.16'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'" + '  -4.67%  '
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'" + '2.186.86'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'" + '  -7.58%  '
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'" + '0.999'
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'" + '  +0.01%  '
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'" + '293.32'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'" + '  -5.44%  '
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'" + '80.59'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'" + '  -7.71%  '
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'" + '  -4.67%  '
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'" + '  +0.06%  '
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'" + '  -7.51%  '
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'" + '0.0769'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'" + '  -8.09%  '
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'" + '27.97'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'" + '  -8.94%  '
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'" + '46.02'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'" + '  -12.36%  '
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'" + '  -2.22%  '
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'" + '2.517.61'
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'" + '  -7.82%  '
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'" + '6.10'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'" + '  -6.61%  '
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'" + '13.80'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'" + '  -7.95%  '
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'" + '2.197.41'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'" + '  -7.52%  '
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'" + '0.706'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'" + '  -6.98%  '
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'" + '38.607.37'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'" + '  -4.70%  '
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'" + '0.0₃0859'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'" + '  -5.37%  '
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'" + '  -8.32%  '
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'" + '64.03'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'" + '  -6.80%  '
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'" + '9.92'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'" + '  -7.92%  '
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'" + '222.70'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'" + '  -4.41%  '
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'" + '  +0.02%  '
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'" + '  -10.26%  '
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'" + '  -4.53%  '
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'" + '22.16'
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'" + '  -6.55%  '
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'" + '  -2.00%  '
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'" + '8.89'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'" + '  -4.64%  '
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'" + '146.65'
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'" + '31.03'
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'" + '  -8.41%  '
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'" + '  +0.00%  '
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'" + '4.75'
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'" + '  -8.78%  '
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'" + '0.0686'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'" + '  -5.93%  '
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'" + '  -5.43%  '
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'" + '0.109'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'" + '  -4.29%  '
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'" + '0.0937'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'" + '  -5.84%  '
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'" + '  -5.88%  '
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'" + '  -7.32%  '
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'" + '14.27'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'" + '  -9.76%  '
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'" + '  -6.94%  '
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'" + '1.880.79'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'" + '  -3.79%  '
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'" + '  -16.61%  '
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'" + '0.0254'
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'" + '  -5.58%  '
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'" + '8.86'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'" + '  -6.86%  '
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'" + '  -10.33%  '
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'" + '2.54'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'" + '  -6.37%  '
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'" + '2.397.45'
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'" + '  -8.01%  '
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'" + '69.74'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'" + '  -3.51%  '
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'" + '85.37'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'" + '  -8.34%  '
$c.Style = "Normal"
